$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "SubUrl"
$ws.Range("E1").Value = "StatuCode"
$ws.Range("D2").Value = "getMultipleServiceListing"
$ws.Range("B2").Value = "ppp8989"
$ws.Range("E2").Value = 200
